$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "next steps" note in F5 to append item 6)
$ws.Range("F5").Value = "My next steps: 1) remove file path  only need file name, 2) see if I can do without the little tk window; 3) output window scroll bar; 4) output window copy-able; 5)put into a stand alone .exe; 6) make a function for parsinng paragraph to sentence"

# Add new log row 6
$ws.Range("B4").Copy($ws.Range("B6"))
$ws.Range("B6").Value = 43839
$ws.Range("D6").Value = $ws.Range("D4").Value2
$ws.Range("E6").Value = "completed #4 and #6 in previous next steps. made parsing paragraph to sentence a function and hide little window that is used to bring up .askdialog prompt."

$ws.Range("C6").Select() | Out-Null
